$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 1 (column E) completed effort for the first five tasks ---
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("E10").Value = 1

# --- New backlog items finished this cycle ---
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = "Criar User Story Primeira Feature"
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = 7
$ws.Range("C12").Value = "Criar User Story Segunda Feature"
$ws.Range("D12").Value = 1

# --- Extra row numbering placeholders ---
$ws.Range("B13").Value = 8
$ws.Range("B14").Value = 9
$ws.Range("B15").Value = 10

# --- Remaining Effort row: now derives from the Completed Effort row instead ---
$ws.Range("E19").Formula = "=D19-E18"
$ws.Range("F19").Formula = "=E19-F18"
$ws.Range("G19").Formula = "=F19-G18"
$ws.Range("H19").Formula = "=G19-H18"
$ws.Range("I19").Formula = "=H19-I18"
$ws.Range("J19").Formula = "=I19-J18"
$ws.Range("K19").Formula = "=J19-K18"

# --- Ideal Burndown row: Day 1 now derived from the Remaining Effort row ---
$ws.Range("E20").Formula = "=D20-E19"

$ws.Application.Calculate()

$ws.Range("F12").Select()
